# Rename the three header/footer logo pictures.
#
# The Pearson logo that appears in the primary footer (id="1") and the
# first-page footer (id="2") is currently named "image2.png"; it should
# become "image1.png". The BTEC logo in the first-page header (id="3")
# is currently named "image1.jpg"; it should become "image2.jpg".
#
# InlineShape has no writable "Name" in the Word object model, so each
# picture is briefly converted to a floating Shape (which does expose
# "Name"), renamed, and converted back to an InlineShape so it stays
# wrapped as <wp:inline> exactly as before.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FooterLogo($footerIndex, $newName) {
    $ftr = $sec.Footers.Item($footerIndex)
    $inl = $ftr.Range.InlineShapes.Item(1)
    $shp = $inl.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

function Rename-HeaderLogo($headerIndex, $newName) {
    $hdr = $sec.Headers.Item($headerIndex)
    $inl = $hdr.Range.InlineShapes.Item(1)
    $shp = $inl.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Primary footer (footer1.xml) - Pearson logo id="1": image2.png -> image1.png
Rename-FooterLogo 1 "image1.png"

# First-page footer (footer2.xml) - Pearson logo id="2": image2.png -> image1.png
Rename-FooterLogo 2 "image1.png"

# First-page header (header2.xml) - BTEC logo id="3": image1.jpg -> image2.jpg
Rename-HeaderLogo 2 "image2.jpg"

Write-Host "Renamed logo pictures in footers/header."
